# Update countries & provincias Spain
# - Refresh "Datos actualizados..." timestamp (10:22 -> 10:52)
# - Refresh COVID stats for several countries; two pairs of countries
#   (Filipinas/Ecuador and Estonia/Marruecos) swapped rank/position,
#   so the country name in column A is updated along with the stats.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp
$ws.Range("A1").Value = "Datos actualizados a 7 de Abril de 2020 a las 10:52"

# Row 17 - Austria: only "Nuevos casos" (F) changes
$ws.Range("F17").Value = 243

# Row 36 - now Filipinas (was Ecuador), fresh stats
$ws.Range("A36").Value = "Filipinas"
$ws.Range("B36").Value = 3764
$ws.Range("C36").Value = 104
$ws.Range("D36").Value = 84
$ws.Range("E36").Value = 3503
$ws.Range("F36").Value = 1
$ws.Range("G36").Value = 14
$ws.Range("H36").Value = 177

# Row 37 - now Ecuador (was Filipinas), keeps Ecuador's previous stats
$ws.Range("A37").Value = "Ecuador"
$ws.Range("B37").Value = 3747
$ws.Range("C37").Value = 0
$ws.Range("D37").Value = 100
$ws.Range("E37").Value = 3456
$ws.Range("F37").Value = 156
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = 191

# Row 61 - now Estonia (was Marruecos), fresh stats
$ws.Range("A61").Value = "Estonia"
$ws.Range("B61").Value = 1149
$ws.Range("C61").Value = 41
$ws.Range("D61").Value = 69
$ws.Range("E61").Value = 1059
$ws.Range("F61").Value = 12
$ws.Range("G61").Value = 2
$ws.Range("H61").Value = 21

# Row 62 - now Marruecos (was Estonia), keeps Marruecos' previous stats
$ws.Range("A62").Value = "Marruecos"
$ws.Range("B62").Value = 1141
$ws.Range("C62").Value = 21
$ws.Range("D62").Value = 88
$ws.Range("E62").Value = 970
$ws.Range("F62").Value = 1
$ws.Range("G62").Value = 3
$ws.Range("H62").Value = 83

# Row 65 - Moldavia: D (Recuperados) and E (Casos activos) change
$ws.Range("D65").Value = 40
$ws.Range("E65").Value = 904

# Row 74 - Kazajistan: B, C, E change
$ws.Range("B74").Value = 685
$ws.Range("C74").Value = 23
$ws.Range("E74").Value = 632

# Row 81 - Letonia: B, C, E change
$ws.Range("B81").Value = 548
$ws.Range("C81").Value = 6
$ws.Range("E81").Value = 531

# Row 101 - Estado de Palestina: B, C, E change
$ws.Range("B101").Value = 260
$ws.Range("C101").Value = 6
$ws.Range("E101").Value = 235
